# Fruta / hortaliza, semanal
# A new weekly record was inserted as row 15, pushing the existing rows
# 15-64 down to 16-65 (the former row 64 is now row 65).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 15, shifting rows 15-64 down to 16-65.
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with the new weekly record.
$ws.Range("A15").Value = 10
$ws.Range("B15").Value = "Vega Modelo de Temuco"
$ws.Range("C15").Value = "La Araucanía"
$ws.Range("D15").Value = 44707
$ws.Range("E15").Value = 9
$ws.Range("F15").Value = 300000001
$ws.Range("G15").Value = "Rabanito"
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 40
$ws.Range("K15").Value = 7000
$ws.Range("L15").Value = 7000
$ws.Range("M15").Value = 7000
$ws.Range("N15").Value = "`$/docena de paquetes"
$ws.Range("O15").Value = "Provincia de Cautín"
$ws.Range("P15").Value = 583
$ws.Range("Q15").Value = 12
$ws.Range("R15").Value = "Hortaliza"
